$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.226637601852417
$ws.Range("B1").Value = 2.365147590637207
$ws.Range("C1").Value = 3.304531335830688
$ws.Range("D1").Value = 3.461189031600952
$ws.Range("E1").Value = 1.117375135421753
